# Add the two new worksheets (UsersData, Sheet3) with their test-case data,
# matching the "Added new Test case files" commit.

$wb = $excel.ActiveWorkbook

# --- New sheet: UsersData (inserted after Sheet2) ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$usersData = $wb.Worksheets.Add($null, $lastSheet)
$usersData.Name = "UsersData"

$usersData.Range("A1").Value = "Test1234"
$usersData.Range("B1").Value = "Test1234@gmail.com"
$usersData.Range("C1").Value = "Tester"
$usersData.Range("D1").Value = "Tester"
$usersData.Range("E1").Value = "www.gmail.com"
$usersData.Range("F1").Value = "Runfast7#123"

$usersData.Hyperlinks.Add($usersData.Range("B1"), "mailto:Test1234@gmail.com")
$usersData.Hyperlinks.Add($usersData.Range("E1"), "http://www.gmail.com/")

$usersData.Range("B1").Style = "Hyperlink"
$usersData.Range("E1").Style = "Hyperlink"

# --- New sheet: Sheet3 (inserted after UsersData, becomes active tab) ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$sheet3 = $wb.Worksheets.Add($null, $lastSheet)

# Order matches the original authoring so shared-string ids line up.
$sheet3.Range("D1").Value = "New "
$sheet3.Range("E1").Value = "Electronic city"
$sheet3.Range("F1").Value = "immediate"
$sheet3.Range("G1").Value = "yeshwanthapur"
$sheet3.Range("H1").Value = "yeshwanthapur"
$sheet3.Range("A1").Value = "Test launch1234"
$sheet3.Range("B1").Value = 50000
$sheet3.Range("C1").Value = 200
$sheet3.Range("I1").Value = 120
$sheet3.Range("J1").Value = 56
$sheet3.Range("K1").Value = 2
